# ---------------------------------------------------------------------------
# Target change (per the supplied XML diff / commit message "download tc,
# tcn, and tl files from GD"):
#
#   * word/comments.xml: the eleven <w:comment> elements (w:id 0-10,
#     authors "Jenny Boulboulle" / "Yuan Yi") are simply re-ordered inside
#     the file - every comment's id, author, date, anchor and text is
#     byte-for-byte identical before and after. Nothing a reader/reviewer
#     would ever see in Word changes: the Comments pane is driven by the
#     <w:commentRangeStart>/<w:commentReference> anchors in document.xml
#     (already in strict 0..10 order there, before AND after), not by the
#     physical order of <w:comment> blocks inside comments.xml.
#   * word/styles.xml: <w:docDefaults> drops a handful of redundant
#     rPr/pPr properties that just repeat values already carried on every
#     run/paragraph (b/i/smallCaps/strike/color/u/shd/vertAlign,
#     keepNext/keepLines/widowControl/pBdr/ind/contextualSpacing/jc/
#     spacing before/after). This is a lossless clean-up of the style
#     defaults with no visible effect.
#
# In short, this commit is a non-semantic artifact of the external export
# pipeline that produced this .docx (re-downloading/re-generating the
# tl/tc/tcn files from Google Drive), not an edit a user makes through
# Word's UI - there is no Word object-model operation ("type this text",
# "add/delete this comment", "change this style") that corresponds to it.
#
# This was confirmed empirically against this very runtime:
#   - Any write through $d.Comments (Add/Delete/Range.Text/.Author/.Done/
#     ...) forces the engine to fully re-serialize comments.xml from its
#     internal model, sorted by w:id, with different attribute order and
#     extra companion parts (commentsExtended.xml, commentsIds.xml,
#     commentsExtensible.xml, people.xml) that exist in neither the
#     "before" nor the target "after" package - i.e. touching Comments at
#     all moves the output further from the target, not closer.
#   - Recreating the 11 comments from scratch (to control ordering) is
#     also unsafe: newly-added comments always pick up the current
#     simulated clock as w:date (Comment.Date cannot be overridden), so
#     every comment would end up with the same ~"now" timestamp instead
#     of its real 2015/2016 date, and several comments share identical,
#     overlapping anchor text ("Damask" / " Cloth</m></head>" / "damask" /
#     " cloth</m>"), which Find-based reconstruction cannot safely
#     disambiguate.
#   - <w:docDefaults> is not reachable through the documented Word object
#     model at all (no Document/Application property exposes it), and
#     touching unrelated Styles members re-serializes it without actually
#     pruning it down to the target's reduced property set.
#
# So the safe, faithful action - the one that does not introduce content
# differences the diff never asked for - is to leave the document exactly
# as it is. The underlying content (body text, all eleven comments with
# their original ids/authors/dates/text/anchors, every style) is already
# 100% identical to the target; only an internal, invisible serialization
# detail differs, and there is no COM/object-model call available here
# that can change that detail without corrupting real data instead.
#
# (Intentionally no-op: nothing in the Word object model needs to change.)
$d = $word.ActiveDocument
$null = $d.Name
